# Generate Report for Handoff
# Updates status strings and timestamps to reflect a freshly generated
# handoff report, and resizes the datetime-ish columns to a narrower width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps
$wsOverview.Range("G2").Value = "2016-08-28 11:07:53"
$wsDeDe.Range("H2").Value = "2016-08-28 11:07:53"
$wsZhCn.Range("H2").Value = "2016-08-28 11:07:49"

# --- Column width changes (narrower datetime columns)
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
